# Apply "bar chart updated (grouped)" edit:
#  - Resumen: recompute Zona/Maximo summary (max over Metricas.Tiempo)
#  - Solucion: re-shuffle the Pedido / Salida assignment table
#  - Metricas: update per-zone Tiempo values

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metricas sheet: updated Tiempo values per Zona
# ---------------------------------------------------------------------
$wsMetricas = $wb.Worksheets.Item("Metricas")
$metricasTiempo = @(700.8299816196345, 454.0177316466646, 676.5996323926912, 491.4706454751866)
for ($i = 0; $i -lt $metricasTiempo.Length; $i++) {
    $wsMetricas.Cells.Item($i + 2, 2).Value = $metricasTiempo[$i]
}

# ---------------------------------------------------------------------
# Resumen sheet: Zona with maximum Tiempo, and that maximum value
# ---------------------------------------------------------------------
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z1"
$wsResumen.Range("C2").Value = 700.8299816196345

# ---------------------------------------------------------------------
# Solucion sheet: Pedido / Salida assignment table (rows 2-81)
# ---------------------------------------------------------------------
$wsSolucion = $wb.Worksheets.Item("Solucion")

$pedidos = @(
    "Pedido_62", "Pedido_14", "Pedido_67", "Pedido_6", "Pedido_29",
    "Pedido_39", "Pedido_7", "Pedido_20", "Pedido_56", "Pedido_37",
    "Pedido_46", "Pedido_28", "Pedido_74", "Pedido_53", "Pedido_23",
    "Pedido_10", "Pedido_12", "Pedido_51", "Pedido_36", "Pedido_59",
    "Pedido_27", "Pedido_65", "Pedido_49", "Pedido_33", "Pedido_18",
    "Pedido_11", "Pedido_52", "Pedido_31", "Pedido_45", "Pedido_61",
    "Pedido_63", "Pedido_21", "Pedido_43", "Pedido_79", "Pedido_77",
    "Pedido_24", "Pedido_13", "Pedido_76", "Pedido_26", "Pedido_9",
    "Pedido_4",  "Pedido_34", "Pedido_50", "Pedido_75", "Pedido_78",
    "Pedido_42", "Pedido_80", "Pedido_47", "Pedido_73", "Pedido_71",
    "Pedido_70", "Pedido_69", "Pedido_64", "Pedido_22", "Pedido_15",
    "Pedido_1",  "Pedido_8",  "Pedido_3",  "Pedido_44", "Pedido_40",
    "Pedido_19", "Pedido_38", "Pedido_32", "Pedido_17", "Pedido_66",
    "Pedido_16", "Pedido_35", "Pedido_55", "Pedido_60", "Pedido_57",
    "Pedido_2",  "Pedido_5",  "Pedido_72", "Pedido_58", "Pedido_25",
    "Pedido_48", "Pedido_30", "Pedido_41", "Pedido_54", "Pedido_68"
)

$salidas = @(
    "S001", "S025", "S041", "S065", "S005",
    "S045", "S029", "S069", "S002", "S042",
    "S026", "S066", "S006", "S030", "S046",
    "S070", "S027", "S003", "S043", "S067",
    "S007", "S031", "S047", "S028", "S004",
    "S071", "S044", "S008", "S068", "S032",
    "S048", "S033", "S009", "S049", "S072",
    "S053", "S037", "S013", "S073", "S050",
    "S010", "S034", "S077", "S054", "S074",
    "S014", "S051", "S038", "S011", "S078",
    "S055", "S015", "S035", "S075", "S012",
    "S052", "S039", "S079", "S056", "S016",
    "S036", "S057", "S017", "S076", "S061",
    "S040", "S058", "S021", "S080", "S018",
    "S062", "S022", "S059", "S063", "S019",
    "S060", "S023", "S064", "S020", "S024"
)

for ($i = 0; $i -lt $pedidos.Length; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $pedidos[$i]
    $wsSolucion.Cells.Item($row, 2).Value = $salidas[$i]
}
